$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows of task/status data appended below the existing table.
$newRows = @(
    @("DNF_1", "done"),
    @("DNF_2", "done"),
    @("DNF_3", "done"),
    @("DNF_4", "done"),
    @("893422_5", "done"),
    @("893422_6", "done"),
    @("DNF_7", "commited")
)

$startRow = 6

# First write all the values for the new rows.
for ($i = 0; $i -lt $newRows.Length; $i++) {
    $r = $startRow + $i
    $task = $newRows[$i][0]
    $status = $newRows[$i][1]

    $ws.Cells.Item($r, 1).Value = $task
    $ws.Cells.Item($r, 2).Value = $status
}

# Task-name column (A) gets a distinct look: vertically centered text...
for ($i = 0; $i -lt $newRows.Length; $i++) {
    $r = $startRow + $i
    $ws.Cells.Item($r, 1).VerticalAlignment = -4108
}

# ...in white, so it stands out against the colored cell background.
for ($i = 0; $i -lt $newRows.Length; $i++) {
    $r = $startRow + $i
    $ws.Cells.Item($r, 1).Font.Color = 0xFFFFFF
}

$ws.Range("B12").Select() | Out-Null
